$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Mapping of (row, col) -> (old, new) text for each arithmetic expression
# that needs updating. Cell ranges are used (rather than a single global
# Find/Replace) because some old values (e.g. "34÷4=") occur more than
# once in the table but must map to different new values.
$edits = @(
    @{ Row = 1;  Col = 1; Old = "42÷9="; New = "97÷8=" },
    @{ Row = 1;  Col = 2; Old = "28÷8="; New = "35÷9=" },
    @{ Row = 1;  Col = 3; Old = "77÷7="; New = "22÷7=" },
    @{ Row = 1;  Col = 4; Old = "42÷6="; New = "38÷5=" },
    @{ Row = 1;  Col = 5; Old = "99÷8="; New = "57÷2=" },

    @{ Row = 5;  Col = 1; Old = "94÷8="; New = "83÷5=" },
    @{ Row = 5;  Col = 2; Old = "34÷4="; New = "49÷5=" },
    @{ Row = 5;  Col = 3; Old = "87÷8="; New = "92÷5=" },
    @{ Row = 5;  Col = 4; Old = "95÷9="; New = "48÷5=" },
    @{ Row = 5;  Col = 5; Old = "39÷6="; New = "93÷5=" },

    @{ Row = 9;  Col = 1; Old = "60÷5="; New = "15÷7=" },
    @{ Row = 9;  Col = 2; Old = "75÷2="; New = "16÷9=" },
    @{ Row = 9;  Col = 3; Old = "74÷7="; New = "40÷4=" },
    @{ Row = 9;  Col = 4; Old = "91÷3="; New = "49÷2=" },
    @{ Row = 9;  Col = 5; Old = "12÷6="; New = "26÷2=" },

    @{ Row = 13; Col = 1; Old = "63÷2="; New = "24÷5=" },
    @{ Row = 13; Col = 2; Old = "59÷3="; New = "34÷2=" },
    @{ Row = 13; Col = 3; Old = "53÷2="; New = "14÷3=" },
    @{ Row = 13; Col = 4; Old = "84÷5="; New = "65÷3=" },
    @{ Row = 13; Col = 5; Old = "18÷7="; New = "15÷3=" },

    @{ Row = 17; Col = 1; Old = "11÷5="; New = "29÷7=" },
    @{ Row = 17; Col = 2; Old = "84÷9="; New = "86÷9=" },
    @{ Row = 17; Col = 3; Old = "34÷5="; New = "43÷9=" },
    @{ Row = 17; Col = 4; Old = "34÷4="; New = "28÷8=" },
    @{ Row = 17; Col = 5; Old = "56÷7="; New = "54÷4=" }
)

foreach ($edit in $edits) {
    $cellRange = $table.Cell($edit.Row, $edit.Col).Range
    # Replace:=1 (wdReplaceOne) so only the single match inside this cell's
    # range is changed, since some old values (e.g. "34÷4=") repeat
    # elsewhere in the table with a different replacement.
    $found = $cellRange.Find.Execute($edit.Old, $true, $false, $false, $false, $false,
                                      $true, 1, $false, $edit.New, 1)
    if (-not $found) {
        Write-Output "WARNING: replacement not found for row=$($edit.Row) col=$($edit.Col) old=$($edit.Old)"
    }
}
